$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1999.6666
$ws.Range("I69").Value = 1999.6666
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 5998.9998
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -5124.9998
$ws.Range("H70").Value = 15264.857
$ws.Range("I70").Value = 700
$ws.Range("K70").Value = 2100
$ws.Range("M70").Value = -1830
$ws.Range("H72").Value = 1999.6666
$ws.Range("I72").Value = 1999.6666
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 17996.9994
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -13628.9994
$ws.Range("H73").Value = 15264.857
$ws.Range("I73").Value = 700
$ws.Range("K73").Value = 2100
$ws.Range("M73").Value = -1164
$ws.Range("H93").Value = 86499
$ws.Range("J93").Value = 86499
$ws.Range("L93").Value = 86499
$ws.Range("N93").Value = -91491
$ws.Range("H99").Value = 1250.4
$ws.Range("I99").Value = 813
$ws.Range("K99").Value = 2439
$ws.Range("M99").Value = -941
$ws.Range("H116").Value = 12537.909
$ws.Range("J116").Value = 3926.4
$ws.Range("L116").Value = 3926.4
$ws.Range("N116").Value = -10810.4
$ws.Range("H118").Value = 836.2
$ws.Range("I118").Value = 836.2
$ws.Range("K118").Value = 2508.6
$ws.Range("M118").Value = -851.6000000000004
$ws.Range("H131").Value = 3577.9
$ws.Range("J131").Value = 4586.4287
$ws.Range("L131").Value = 13759.2861
$ws.Range("N131").Value = -23839.2861
$ws.Range("H137").Value = 1435.3636
$ws.Range("I137").Value = 1361.125
$ws.Range("J137").Value = 1633.3334
$ws.Range("K137").Value = 4083.375
$ws.Range("L137").Value = 4900.0002
$ws.Range("M137").Value = -1533.375
$ws.Range("N137").Value = -10000.0002
$ws.Range("H138").Value = 3188.9148
$ws.Range("I138").Value = 3587.5789
$ws.Range("J138").Value = 2918.3928
$ws.Range("K138").Value = 10762.7367
$ws.Range("L138").Value = 8755.178400000001
$ws.Range("M138").Value = -5622.736699999999
$ws.Range("N138").Value = -19035.1784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4937.41
$ws.Range("I32").Value = 3631.0322
$ws.Range("K32").Value = 3631.0322
$ws.Range("M32").Value = -3344.0322
$ws.Range("H61").Value = 4345.6816
$ws.Range("J61").Value = 9699.799999999999
$ws.Range("L61").Value = 9699.799999999999
$ws.Range("N61").Value = -10123.8
$ws.Range("H97").Value = 535.2381
$ws.Range("I97").Value = 496.30768
$ws.Range("J97").Value = 598.5
$ws.Range("K97").Value = 496.30768
$ws.Range("L97").Value = 598.5
$ws.Range("M97").Value = -0.3076800000000048
$ws.Range("N97").Value = -1590.5
$ws.Range("H132").Value = 1616.6428
$ws.Range("I132").Value = 1012.05884
$ws.Range("J132").Value = 2551
$ws.Range("K132").Value = 3036.17652
$ws.Range("L132").Value = 7653
$ws.Range("M132").Value = -506.17652
$ws.Range("N132").Value = -12713
$ws.Range("H136").Value = 4345.6816
$ws.Range("J136").Value = 9699.799999999999
$ws.Range("L136").Value = 29099.4
$ws.Range("N136").Value = -34199.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 315.25806
$ws.Range("I94").Value = 320.76666
$ws.Range("K94").Value = 320.76666
$ws.Range("M94").Value = 130.23334
$ws.Range("H105").Value = 1984.9546
$ws.Range("I105").Value = 1962.1666
$ws.Range("J105").Value = 2087.5
$ws.Range("K105").Value = 1962.1666
$ws.Range("L105").Value = 2087.5
$ws.Range("M105").Value = -215.1666
$ws.Range("N105").Value = -5581.5
$ws.Range("H134").Value = 13523.782
$ws.Range("I134").Value = 13141.556
$ws.Range("K134").Value = 39424.66800000001
$ws.Range("M134").Value = -36889.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 734.5454999999999
$ws.Range("I16").Value = 720.2222
$ws.Range("J16").Value = 799
$ws.Range("K16").Value = 720.2222
$ws.Range("L16").Value = 799
$ws.Range("M16").Value = -433.2222
$ws.Range("N16").Value = -1373
$ws.Range("H31").Value = 3120.7896
$ws.Range("I31").Value = 2514.8147
$ws.Range("K31").Value = 2514.8147
$ws.Range("M31").Value = -2219.8147
$ws.Range("H34").Value = 3120.7896
$ws.Range("I34").Value = 2514.8147
$ws.Range("K34").Value = 2514.8147
$ws.Range("M34").Value = -2312.8147
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1706
$ws.Range("H113").Value = 734.5454999999999
$ws.Range("I113").Value = 720.2222
$ws.Range("J113").Value = 799
$ws.Range("K113").Value = 720.2222
$ws.Range("L113").Value = 799
$ws.Range("M113").Value = 1449.7778
$ws.Range("N113").Value = -5139
$ws.Range("H122").Value = 2315.5
$ws.Range("I122").Value = 1841.9166
$ws.Range("J122").Value = 5157
$ws.Range("K122").Value = 5525.7498
$ws.Range("L122").Value = 15471
$ws.Range("M122").Value = -3075.7498
$ws.Range("N122").Value = -20371
$ws.Range("H134").Value = 995.88
$ws.Range("I134").Value = 857.0454999999999
$ws.Range("K134").Value = 2571.1365
$ws.Range("M134").Value = -36.13649999999961

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2973.3333
$ws.Range("I80").Value = 2872.2856
$ws.Range("K80").Value = 2872.2856
$ws.Range("M80").Value = -1874.2856
$ws.Range("H83").Value = 2973.3333
$ws.Range("I83").Value = 2872.2856
$ws.Range("K83").Value = 14361.428
$ws.Range("M83").Value = -9369.428
$ws.Range("H102").Value = 3205.182
$ws.Range("I102").Value = 3417.3333
$ws.Range("K102").Value = 3417.3333
$ws.Range("M102").Value = -1795.3333
$ws.Range("H107").Value = 1367.3334
$ws.Range("I107").Value = 999
$ws.Range("K107").Value = 999
$ws.Range("M107").Value = 921
$ws.Range("H122").Value = 1954.5385
$ws.Range("I122").Value = 1400.3334
$ws.Range("J122").Value = 2429.5715
$ws.Range("K122").Value = 4201.0002
$ws.Range("L122").Value = 7288.7145
$ws.Range("M122").Value = -1751.0002
$ws.Range("N122").Value = -12188.7145
$ws.Range("H126").Value = 34828.645
$ws.Range("I126").Value = 2775.0454
$ws.Range("K126").Value = 8325.136200000001
$ws.Range("M126").Value = -5855.136200000001
$ws.Range("H132").Value = 3290.72
$ws.Range("I132").Value = 2643.55
$ws.Range("K132").Value = 7930.650000000001
$ws.Range("M132").Value = -5400.650000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 397500
$ws.Range("H7").Value = 5698.778
$ws.Range("I7").Value = 3280.4
$ws.Range("J7").Value = 8721.75
$ws.Range("K7").Value = 3280.4
$ws.Range("L7").Value = 8721.75
$ws.Range("M7").Value = -3168.4
$ws.Range("N7").Value = -8945.75
$ws.Range("H40").Value = 10666.333
$ws.Range("I40").Value = 1999.5
$ws.Range("J40").Value = 14999.75
$ws.Range("K40").Value = 1999.5
$ws.Range("L40").Value = 14999.75
$ws.Range("M40").Value = -1863.5
$ws.Range("N40").Value = -15271.75
$ws.Range("H46").Value = 1532.125
$ws.Range("I46").Value = 864.75
$ws.Range("K46").Value = 864.75
$ws.Range("M46").Value = -676.75
$ws.Range("H68").Value = 2789.7273
$ws.Range("I68").Value = 2520.7778
$ws.Range("K68").Value = 2520.7778
$ws.Range("M68").Value = -1771.7778
$ws.Range("H71").Value = 2789.7273
$ws.Range("I71").Value = 2520.7778
$ws.Range("K71").Value = 12603.889
$ws.Range("M71").Value = -8859.888999999999
$ws.Range("H93").Value = 376.09525
$ws.Range("I93").Value = 332.84616
$ws.Range("J93").Value = 446.375
$ws.Range("K93").Value = 332.84616
$ws.Range("L93").Value = 446.375
$ws.Range("M93").Value = 915.1538399999999
$ws.Range("N93").Value = -2942.375
$ws.Range("H95").Value = 94000
$ws.Range("J95").Value = 94000
$ws.Range("L95").Value = 94000
$ws.Range("N95").Value = -99492
$ws.Range("H122").Value = 9800.799999999999
$ws.Range("I122").Value = 8501
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 25503
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -23053
$ws.Range("N122").Value = -49900
$ws.Range("H126").Value = 5698.778
$ws.Range("I126").Value = 3280.4
$ws.Range("J126").Value = 8721.75
$ws.Range("K126").Value = 9841.200000000001
$ws.Range("L126").Value = 26165.25
$ws.Range("M126").Value = -7371.200000000001
$ws.Range("N126").Value = -31105.25
$ws.Range("H132").Value = 1995.1305
$ws.Range("I132").Value = 1775.4546
$ws.Range("J132").Value = 2196.5
$ws.Range("K132").Value = 5326.3638
$ws.Range("L132").Value = 6589.5
$ws.Range("M132").Value = -2796.3638
$ws.Range("N132").Value = -11649.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 31844
$ws.Range("I51").Value = 24800
$ws.Range("J51").Value = 38888
$ws.Range("K51").Value = 24800
$ws.Range("L51").Value = 38888
$ws.Range("M51").Value = -24290
$ws.Range("N51").Value = -39908
$ws.Range("H81").Value = 1393.6666
$ws.Range("I81").Value = 1138.2222
$ws.Range("K81").Value = 2276.4444
$ws.Range("M81").Value = -1215.4444
$ws.Range("H84").Value = 1393.6666
$ws.Range("I84").Value = 1138.2222
$ws.Range("K84").Value = 11382.222
$ws.Range("M84").Value = -6078.222
$ws.Range("H136").Value = 2735.6897
$ws.Range("I136").Value = 3011.5
$ws.Range("J136").Value = 2478.2666
$ws.Range("K136").Value = 9034.5
$ws.Range("L136").Value = 7434.7998
$ws.Range("M136").Value = -6484.5
$ws.Range("N136").Value = -12534.7998
